$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark these BOM lines as fully "Have" (column C = Qty from column B),
# which drives the "Need" formula in column E back down to 0.
$ws.Range("C2").Value  = 1
$ws.Range("C3").Value  = 4
$ws.Range("C6").Value  = 2
$ws.Range("C7").Value  = 6
$ws.Range("C8").Value  = 2
$ws.Range("C9").Value  = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 1

# J109's connector part was actually the 10-pin Eurorack power connector,
# not the 16-pin one.
$ws.Range("F15").Value = "Conn_Eurorack_10"

# Leave the selection where the author ended up after editing.
$ws.Range("C16").Select() | Out-Null
